$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 17 (shifts existing rows 17..59 down to 18..60)
$ws.Rows.Item(17).Insert()

# Populate the newly inserted row 17 with the new weekly record
$ws.Range("A17").Value = 11
$ws.Range("B17").Value = "Vega Monumental Concepción"
$ws.Range("C17").Value = "Bíobío"
$ws.Range("D17").Value = 44565
$ws.Range("E17").Value = 8
$ws.Range("F17").Value = 100112001
$ws.Range("G17").Value = "Berenjena"
$ws.Range("H17").Value = "Sin especificar"
$ws.Range("I17").Value = "Primera"
$ws.Range("J17").Value = 100
$ws.Range("K17").Value = 11000
$ws.Range("L17").Value = 12000
$ws.Range("M17").Value = 11500
$ws.Range("N17").Value = "$/caja 60 unidades"
$ws.Range("O17").Value = "Región Metropolitana"
$ws.Range("P17").Value = 192
$ws.Range("Q17").Value = 60
$ws.Range("R17").Value = "Hortaliza"
